$wb = $excel.ActiveWorkbook

function Set-SheetRows {
    # Positional params only — named parameter binding is unreliable here.
    param($ws, $oldLastRow, $newLastRow, $lastCol, $rows)

    # Clear the existing data block (keep header row 1 intact).
    if ($oldLastRow -ge 2) {
        $clearRange = $ws.Cells.Item(2, 1).Resize($oldLastRow - 1, $lastCol)
        $clearRange.ClearContents()
    }

    # Grow / shrink the sheet so the used range ends up exactly at newLastRow.
    if ($newLastRow -lt $oldLastRow) {
        $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 1)).EntireRow.Delete()
    }
    elseif ($newLastRow -gt $oldLastRow) {
        $ws.Range($ws.Cells.Item($oldLastRow + 1, 1), $ws.Cells.Item($newLastRow, 1)).EntireRow.Insert()
    }

    # Write the new values.
    $r = 2
    foreach ($row in $rows) {
        $c = 1
        foreach ($val in $row) {
            $ws.Cells.Item($r, $c).Value = $val
            $c++
        }
        $r++
    }
}

# ---------------------------------------------------------------------------
# Sheet: Home win  (A1:F6 -> A1:F5)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Home win")
$rows1 = @(
    , @("17-01-2025 19:00", "MOROCCO", "BOTOLA PRO", "CODM Meknès - Moghreb Tetouan", 80, 2.1)
    , @("17-01-2025 19:45", "NORTHERN-IRELAND", "PREMIERSHIP", "Dungannon Swifts - Glenavon FC", 73.3, 1.95)
    , @("17-01-2025 17:30", "SOUTH-AFRICA", "PREMIER SOCCER LEAGUE", "Cape Town City - Richards Bay", 73.3, 1.95)
    , @("17-01-2025 19:30", "SPAIN", "PRIMERA DIVISIÓN RFEF - GROUP 1", "Sestao River - Real Unión", 80, 1.95)
)
Set-SheetRows $ws1 6 5 6 $rows1

# ---------------------------------------------------------------------------
# Sheet: Draw  (A1:F3 -> A1:F2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Draw")
$rows2 = @(
    , @("17-01-2025 19:00", "FRANCE", "LIGUE 2", "Guingamp - Rodez", 60, 3.5)
)
Set-SheetRows $ws2 3 2 6 $rows2

# ---------------------------------------------------------------------------
# Sheet: Btts  (A1:F9 -> A1:F7)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Btts")
$rows3 = @(
    , @("17-01-2025 16:00", "BAHRAIN", "PREMIER LEAGUE", "Bahrain SC - Sitra", 93.3, 1.8)
    , @("17-01-2025 18:30", "FRANCE", "NATIONAL 1", "Orleans - Chateauroux", 80, 2.05)
    , @("17-01-2025 08:30", "INDONESIA", "LIGA 1", "Persebaya Surabaya - Malut United", 96, 1.9)
    , @("17-01-2025 19:30", "ITALY", "SERIE C - GIRONE C", "Giugliano - Audace Cerignola", 81.7, 1.8)
    , @("17-01-2025 19:30", "SPAIN", "PRIMERA DIVISIÓN RFEF - GROUP 2", "Alcorcon - Fuenlabrada", 86.7, 1.95)
    , @("17-01-2025 17:00", "TURKEY", "1. LIG", "Şanlıurfaspor - Ankaragucu", 76, 1.8)
)
Set-SheetRows $ws3 9 7 6 $rows3

# ---------------------------------------------------------------------------
# Sheet: Over_Under  (A1:H6 -> A1:H8)
#   Rows 2-4 unchanged; row5 H-odds change; a new row6 is inserted; the old
#   row6 (Hannover) shifts to row7 with updated odds; a new row8 is appended.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Over_Under")

# H3: 2.3 -> 2.2
$ws4.Cells.Item(3, 8).Value = 2.2

# H5: 2.38 -> 2.5
$ws4.Cells.Item(5, 8).Value = 2.5

# Insert a new row before the current row 6 (shifts old row6 -> row7).
$ws4.Range($ws4.Cells.Item(6, 1), $ws4.Cells.Item(6, 8)).EntireRow.Insert()

# New row6 content.
$newRow6 = @("17-01-2025 19:00", "BELGIUM", "CHALLENGER PRO LEAGUE", "Patro Eisden - Lokeren-Temse", 60, 1.87, 60, 3.1)
$c = 1
foreach ($val in $newRow6) {
    $ws4.Cells.Item(6, $c).Value = $val
    $c++
}

# Old Hannover row is now row7; update its odds (F: 1.53 -> 1.55, H: 2.3 -> 2.4).
$ws4.Cells.Item(7, 6).Value = 1.55
$ws4.Cells.Item(7, 8).Value = 2.4

# Append the new row8.
$newRow8 = @("17-01-2025 10:30", "ISRAEL", "LIGA ALEF", "Tzeirei Umm Al-Fahm - Kiryat Yam SC", 70, 1.61, 65, 2.5)
$c = 1
foreach ($val in $newRow8) {
    $ws4.Cells.Item(8, $c).Value = $val
    $c++
}

# ---------------------------------------------------------------------------
# Sheet: Away Win  (A1:F2 unchanged size, content replaced)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Away Win")
$rows5 = @(
    , @("17-01-2025 11:00", "ISRAEL", "LIGA ALEF", "Maccabi Kiryat Malachi - Maccabi Ashdod", 73.3, 2.45)
)
Set-SheetRows $ws5 2 2 6 $rows5
